# Update the "想去人数" (F column) counts for several rows across the
# 展览, 演出 and 全部类型 worksheets, as generated by the latest scrape
# (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map: sheet name -> list of (row, new value)
$updates = @{
    "展览" = @(
        @{ Row = 4;  Value = 9567 },
        @{ Row = 13; Value = 453 },
        @{ Row = 14; Value = 12192 },
        @{ Row = 33; Value = 595 },
        @{ Row = 35; Value = 3063 },
        @{ Row = 40; Value = 24 },
        @{ Row = 41; Value = 112 },
        @{ Row = 42; Value = 432 },
        @{ Row = 43; Value = 548 },
        @{ Row = 45; Value = 139 },
        @{ Row = 47; Value = 114 }
    )
    "演出" = @(
        @{ Row = 12; Value = 50 },
        @{ Row = 26; Value = 37 }
    )
    "全部类型" = @(
        @{ Row = 5;  Value = 9567 },
        @{ Row = 14; Value = 453 },
        @{ Row = 15; Value = 12192 },
        @{ Row = 21; Value = 50 },
        @{ Row = 26; Value = 2102 },
        @{ Row = 34; Value = 595 },
        @{ Row = 36; Value = 3063 },
        @{ Row = 40; Value = 112 },
        @{ Row = 41; Value = 432 },
        @{ Row = 42; Value = 37 },
        @{ Row = 43; Value = 548 },
        @{ Row = 45; Value = 139 },
        @{ Row = 47; Value = 114 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
